$p = $ppt.ActivePresentation

# The commit duplicates slide 2 ("BEFORE WE GET STARTED") and the
# duplicate is inserted immediately after it, becoming the new slide 3.
# All the slides that used to follow (USING GIT & GITHUB REPOSITORY,
# GIT BASH & COMMAND LINE, TRY GITHUB REPOSITORY, USING SOME RESOURCE,
# LET'S DO IT TOGETHER !!!) simply shift down by one position as a
# result - no other content changes are required.
$s2 = $p.Slides.Item(2)
$s2.Duplicate() | Out-Null
